# VGG10 results.xlsx update
# - LeNet9 for CIFAR10 and VGG10 for MNIST: fill in the measured Acck / Acc
#   (columns D, E, H) for the first two "Base Model: VGG10 / TASK: MNIST"
#   blocks (rows 3-10, MLP encoder/decoder; rows 12-25, Conv encoder/decoder).
#   Columns I (Acc desc), J (RA0) and K (RAk) are formulas and recompute
#   automatically once D/E/H are populated.
# - Restore the sheet/window view state that was active when the file was
#   last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Block 1: rows 3-6 (K=2,R=1,N=3 / "train")
# ---------------------------------------------------------------------
$ws.Range("D3:D6").Value = 0.99399999999999999
$ws.Range("E3:E6").Value = 0.99350000000000005

$ws.Range("H3").Value = 0.99350000000000005
$ws.Range("H4").Value = 0.94089999999999996
$ws.Range("H5").Value = 0.63049999999999995
$ws.Range("H6").Value = 0.1124

# ---------------------------------------------------------------------
# Block 2: rows 7-10 (K=2,R=1,N=3 / "test")
# ---------------------------------------------------------------------
$ws.Range("D7:D10").Value = 0.99139999999999995
$ws.Range("E7:E10").Value = 0.99050000000000005

$ws.Range("H7").Value = 0.99050000000000005
$ws.Range("H8").Value = 0.93110000000000004
$ws.Range("H9").Value = 0.7
$ws.Range("H10").Value = 0.1135

# ---------------------------------------------------------------------
# Block 3: rows 12-18 (K=4,R=2,N=6 / "train")
# ---------------------------------------------------------------------
$ws.Range("D12:D18").Value = 0.99399999999999999
$ws.Range("E12:E18").Value = 0.9929

$ws.Range("H12").Value = 0.9929
$ws.Range("H13").Value = 0.98740000000000006
$ws.Range("H14").Value = 0.97209999999999996
$ws.Range("H15").Value = 0.9244
$ws.Range("H16").Value = 0.83630000000000004
$ws.Range("H17").Value = 0.61419999999999997
$ws.Range("H18").Value = 0.1124

# ---------------------------------------------------------------------
# Block 4: rows 19-25 (K=4,R=2,N=6 / "test")
# ---------------------------------------------------------------------
$ws.Range("D19:D25").Value = 0.99139999999999995
$ws.Range("E19:E25").Value = 0.9899

$ws.Range("H19").Value = 0.9899
$ws.Range("H20").Value = 0.98670000000000002
$ws.Range("H21").Value = 0.96909999999999996
$ws.Range("H22").Value = 0.9304
$ws.Range("H23").Value = 0.83689999999999998
$ws.Range("H24").Value = 0.60560000000000003
$ws.Range("H25").Value = 0.1135

# ---------------------------------------------------------------------
# Restore view state: active window scrolled to D16, selection on G25.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$ws.Range("D16").Select()
$win.ScrollRow = 16
$win.ScrollColumn = 4
$ws.Range("G25").Select()

$wb.Saved = $false
